# "Register and signin Link" — clear out the stale/sample login test data
# from the "login" sheet, leaving only the header row and a single sample
# row (Testing12 / Test123), and move the active selection to C7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

# Row 2 used to hold a leftover username sample ("Testing1") — clear it.
$ws.Range("A2:D2").ClearContents()

# Row 3 (Testing12 / Test123) stays as-is; only make sure the trailing
# columns remain empty.
$ws.Range("C3:D3").ClearContents()

# Rows 4-9 held a grid of sample username/password/confirmation data that
# is no longer needed — wipe all of it.
$ws.Range("A4:D9").ClearContents()

# Move the active cell/selection to C7, matching the saved view state.
$ws.Activate()
$ws.Range("C7").Select()
